$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.579.10'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.839.92'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5273'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3144'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06793'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7787'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07754'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '1.854.76'
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.005'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9999'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9997'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007905'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").Value = '26.595.46'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Value = '2.078.08'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.597'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.313'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.10%  '
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.207'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.684'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '110.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.167'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.062'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04865'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7301'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.858'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.093'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.239'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01719'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4798'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8943'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.913'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.639'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4151'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.032'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1234'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05816'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.94%  '
